$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 130992753
$ws.Range("M7").Value = "färska spår"
$ws.Range("P7").Value = "Skansen, Skansen, Boh"
$ws.Range("Z7").Value = "14:36"
$ws.Range("AB7").Value = "14:36"
